$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17 (idx 0)
$ws.Range("H17").Value = 742957.4399999999
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 742957.4399999999
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2228872.32
$ws.Range("N17").Value = -2229208.32
# row 24 (idx 1)
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = ""
# row 32 (idx 2)
$ws.Range("H32").Value = 2555
$ws.Range("I32").Value = 2555
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2555
$ws.Range("L32").Value = 0
# row 38 (idx 3)
$ws.Range("H38").Value = 77.916664
$ws.Range("I38").Value = 83.181816
$ws.Range("J38").Value = 20
$ws.Range("K38").Value = 249.545448
$ws.Range("L38").Value = 60
$ws.Range("M38").Value = 122.454552
# row 39 (idx 4)
$ws.Range("H39").Value = 318.8889
$ws.Range("I39").Value = 40
$ws.Range("J39").Value = 458.33334
$ws.Range("K39").Value = 120
$ws.Range("L39").Value = 1375.00002
$ws.Range("M39").Value = 176
$ws.Range("N39").Value = -1967.00002
# row 42 (idx 5)
$ws.Range("H42").Value = 239.8
$ws.Range("I42").Value = 67.666664
$ws.Range("J42").Value = 498
$ws.Range("K42").Value = 202.999992
$ws.Range("L42").Value = 1494
$ws.Range("M42").Value = 27.00000800000001
$ws.Range("N42").Value = -1954
# row 47 (idx 6)
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
# row 96 (idx 7)
$ws.Range("H96").Value = 2196.1428
$ws.Range("I96").Value = 2228.8333
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 6686.499899999999
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -5313.499899999999
$ws.Range("N96").Value = -8746
# row 97 (idx 8)
$ws.Range("H97").Value = 18968.111
$ws.Range("I97").Value = 925
$ws.Range("J97").Value = 24123.285
$ws.Range("K97").Value = 2775
$ws.Range("L97").Value = 72369.855
$ws.Range("M97").Value = -2279
$ws.Range("N97").Value = -73361.855
# row 101 (idx 9)
$ws.Range("H101").Value = 297.2
$ws.Range("I101").Value = 331
$ws.Range("J101").Value = 246.5
$ws.Range("K101").Value = 993
$ws.Range("L101").Value = 739.5
$ws.Range("M101").Value = 629
$ws.Range("N101").Value = -3983.5
# row 115 (idx 10)
$ws.Range("H115").Value = 218.5
$ws.Range("I115").Value = 218.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 655.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = ""
# row 132 (idx 11)
$ws.Range("H132").Value = 2052.8276
$ws.Range("I132").Value = 1442.4546
$ws.Range("J132").Value = 3971.1428
$ws.Range("K132").Value = 4327.3638
$ws.Range("L132").Value = 11913.4284
$ws.Range("M132").Value = -1797.3638
$ws.Range("N132").Value = -16973.4284

$ws = $wb.Worksheets.Item("ARM")
# row 2 (idx 12)
$ws.Range("H2").Value = 1957.4445
$ws.Range("I2").Value = 2136.8125
$ws.Range("J2").Value = 522.5
$ws.Range("K2").Value = 2136.8125
$ws.Range("L2").Value = 522.5
$ws.Range("M2").Value = -2023.8125
$ws.Range("N2").Value = -748.5
# row 27 (idx 13)
$ws.Range("H27").Value = 7498
$ws.Range("I27").Value = 7498
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 7498
$ws.Range("L27").Value = ""
# row 45 (idx 14)
$ws.Range("H45").Value = 92071.766
$ws.Range("I45").Value = 116449.5
$ws.Range("J45").Value = 10812.667
$ws.Range("K45").Value = 116449.5
$ws.Range("L45").Value = 10812.667
$ws.Range("M45").Value = -116072.5
$ws.Range("N45").Value = -11566.667
# row 74 (idx 15)
$ws.Range("H74").Value = 7954.607
$ws.Range("I74").Value = 1251.3529
$ws.Range("J74").Value = 18314.182
$ws.Range("K74").Value = 1251.3529
$ws.Range("L74").Value = 18314.182
$ws.Range("M74").Value = -377.3529000000001
$ws.Range("N74").Value = -20062.182
# row 77 (idx 16)
$ws.Range("H77").Value = 7954.607
$ws.Range("I77").Value = 1251.3529
$ws.Range("J77").Value = 18314.182
$ws.Range("K77").Value = 6256.7645
$ws.Range("L77").Value = 91570.91
$ws.Range("M77").Value = -1888.7645
$ws.Range("N77").Value = -100306.91
# row 110 (idx 17)
$ws.Range("H110").Value = 7029.7085
$ws.Range("I110").Value = 7372.9443
$ws.Range("J110").Value = 6000
$ws.Range("K110").Value = 7372.9443
$ws.Range("L110").Value = 6000
$ws.Range("M110").Value = -5327.9443
$ws.Range("N110").Value = -10090
# row 116 (idx 18)
$ws.Range("H116").Value = 1957.4445
$ws.Range("I116").Value = 2136.8125
$ws.Range("J116").Value = 522.5
$ws.Range("K116").Value = 2136.8125
$ws.Range("L116").Value = 522.5
$ws.Range("M116").Value = 157.1875
$ws.Range("N116").Value = -5110.5

$ws = $wb.Worksheets.Item("BSM")
# row 3 (idx 19)
$ws.Range("H3").Value = 1957.4445
$ws.Range("I3").Value = 2136.8125
$ws.Range("J3").Value = 522.5
$ws.Range("K3").Value = 2136.8125
$ws.Range("L3").Value = 522.5
$ws.Range("M3").Value = -2022.8125
$ws.Range("N3").Value = -750.5
# row 20 (idx 20)
$ws.Range("H20").Value = 3872.1667
$ws.Range("I20").Value = 3192.5334
$ws.Range("J20").Value = 4551.8
$ws.Range("K20").Value = 3192.5334
$ws.Range("L20").Value = 4551.8
$ws.Range("M20").Value = -2945.5334
$ws.Range("N20").Value = -5045.8
# row 86 (idx 21)
$ws.Range("H86").Value = 2018.75
$ws.Range("I86").Value = 2018.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2018.75
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -895.75
# row 89 (idx 22)
$ws.Range("H89").Value = 2018.75
$ws.Range("I89").Value = 2018.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10093.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4477.75
# row 105 (idx 23)
$ws.Range("H105").Value = 2011.762
$ws.Range("I105").Value = 1486
$ws.Range("J105").Value = 5166.3335
$ws.Range("K105").Value = 1486
$ws.Range("L105").Value = 5166.3335
$ws.Range("M105").Value = 261
$ws.Range("N105").Value = -8660.333500000001
# row 107 (idx 24)
$ws.Range("H107").Value = 1881.8572
$ws.Range("I107").Value = 1447.6
$ws.Range("J107").Value = 2967.5
$ws.Range("K107").Value = 1447.6
$ws.Range("L107").Value = 2967.5
$ws.Range("M107").Value = 472.4000000000001
$ws.Range("N107").Value = -6807.5
# row 132 (idx 25)
$ws.Range("H132").Value = 100589.8
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 100589.8
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 100589.8
$ws.Range("N132").Value = -110709.8
# row 134 (idx 26)
$ws.Range("H134").Value = 1727.4615
$ws.Range("I134").Value = 1716.56
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 5149.68
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -2614.68
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
# row 7 (idx 27)
$ws.Range("H7").Value = 90.44444
$ws.Range("I7").Value = 62.4
$ws.Range("J7").Value = 125.5
$ws.Range("K7").Value = 62.4
$ws.Range("L7").Value = 125.5
$ws.Range("M7").Value = 50.6
$ws.Range("N7").Value = -351.5
# row 13 (idx 28)
$ws.Range("H13").Value = 2460.6
$ws.Range("I13").Value = 1601.3334
$ws.Range("J13").Value = 3749.5
$ws.Range("K13").Value = 1601.3334
$ws.Range("L13").Value = 3749.5
$ws.Range("M13").Value = -1462.3334
$ws.Range("N13").Value = -4027.5
# row 92 (idx 29)
$ws.Range("H92").Value = 18550.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 18550.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 18550.5

$ws = $wb.Worksheets.Item("CUL")
# row 12 (idx 30)
$ws.Range("H12").Value = 458.54544
$ws.Range("I12").Value = 6.3333335
$ws.Range("J12").Value = 628.125
$ws.Range("K12").Value = 19.0000005
$ws.Range("L12").Value = 1884.375
$ws.Range("M12").Value = 153.9999995
$ws.Range("N12").Value = -2230.375

$ws = $wb.Worksheets.Item("GSM")
# row 53 (idx 31)
$ws.Range("H53").Value = 19999.5
$ws.Range("I53").Value = 19999.5
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 19999.5
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -19368.5
# row 134 (idx 32)
$ws.Range("H134").Value = 46150
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 46150
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 138450
$ws.Range("N134").Value = -143520

$ws = $wb.Worksheets.Item("LTW")
# row 16 (idx 33)
$ws.Range("H16").Value = 853.2727
$ws.Range("I16").Value = 773.3182
$ws.Range("J16").Value = 1013.1818
$ws.Range("K16").Value = 773.3182
$ws.Range("L16").Value = 1013.1818
$ws.Range("M16").Value = -603.3182
$ws.Range("N16").Value = -1353.1818
# row 22 (idx 34)
$ws.Range("H22").Value = 2099.2334
$ws.Range("I22").Value = 1749
$ws.Range("J22").Value = 2332.7222
$ws.Range("K22").Value = 1749
$ws.Range("L22").Value = 2332.7222
$ws.Range("M22").Value = -1454
$ws.Range("N22").Value = -2922.7222
# row 27 (idx 35)
$ws.Range("H27").Value = 2099.2334
$ws.Range("I27").Value = 1749
$ws.Range("J27").Value = 2332.7222
$ws.Range("K27").Value = 1749
$ws.Range("L27").Value = 2332.7222
$ws.Range("M27").Value = -1642
$ws.Range("N27").Value = -2546.7222
# row 93 (idx 36)
$ws.Range("H93").Value = 333866.66
$ws.Range("I93").Value = 500500
$ws.Range("J93").Value = 600
$ws.Range("K93").Value = 500500
$ws.Range("L93").Value = 600
$ws.Range("M93").Value = -499252
$ws.Range("N93").Value = -3096
# row 132 (idx 37)
$ws.Range("H132").Value = 4287.1904
$ws.Range("I132").Value = 4207.4116
$ws.Range("J132").Value = 4626.25
$ws.Range("K132").Value = 12622.2348
$ws.Range("L132").Value = 13878.75
$ws.Range("M132").Value = -10092.2348
$ws.Range("N132").Value = -18938.75

$ws = $wb.Worksheets.Item("WVR")
# row 74 (idx 38)
$ws.Range("H74").Value = 16424.6
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 16424.6
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 16424.6
$ws.Range("N74").Value = -18296.6
# row 77 (idx 39)
$ws.Range("H77").Value = 16424.6
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 16424.6
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 49273.8
$ws.Range("N77").Value = -58633.8
# row 100 (idx 40)
$ws.Range("H100").Value = 1711.3077
$ws.Range("I100").Value = 800
$ws.Range("J100").Value = 1787.25
$ws.Range("K100").Value = 1600
$ws.Range("L100").Value = 3574.5
$ws.Range("M100").Value = -1059
$ws.Range("N100").Value = -4656.5
# row 107 (idx 41)
$ws.Range("H107").Value = 569.5714
$ws.Range("I107").Value = 515.75
$ws.Range("J107").Value = 641.3333
$ws.Range("K107").Value = 1547.25
$ws.Range("L107").Value = 1923.9999
$ws.Range("M107").Value = 372.75
$ws.Range("N107").Value = -5763.9999
